# PARCIAL II / PRIMER PUNTO.xlsx - apply author's edit
# - Label the a=/b= (Y range) and c=/d= (X range) rows with "Y"/"X" in column C
# - Draw a full grid border ("All Borders") over both mini result tables
#   (A15:D41 and F15:I21), reusing/creating the bold+center header styles
# - Move the visible selection to D8 (and let the view's top-left reset)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Label the range-bounds rows (a=, b= -> "Y"; c=, d= -> "X") ---
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"
$ws.Range("C8").Value = "X"
$ws.Range("C9").Value = "X"

# --- Apply "All Borders" (thin, every side of every cell) to both tables ---
$leftTable = $ws.Range("A15:D41")
$leftTable.Borders.LineStyle = 1

$rightTable = $ws.Range("F15:I21")
$rightTable.Borders.LineStyle = 1

# --- Update selection / scroll position ---
[void]$ws.Range("D8").Select()
